# Generate Report for Archive
#
# The CI run that produced this report picked up the f843c243... file's
# status flipping to "In Translation" (it had been re-queued), which moves
# it up in the localization-status report so the three in-flight files are
# grouped: f843c243 (now row 7, still In Translation), then 5030e9b0 (row 8)
# and b9bdd042 (row 9), both still "Ready for handoff".
#
# Overview sheet (File Name / zh-cn / de-de / Latest Handoff Date)
$wsOverview = $excel.ActiveWorkbook.Worksheets.Item("Overview")

$wsOverview.Range("A7").Value = "f843c243-6809-434c-80ed-ff15e90cb237.md"
$wsOverview.Range("B7").Value = "In Translation"
$wsOverview.Range("C7").Value = "In Translation"
$wsOverview.Range("D7").Value = "2016-30-12 12:30:56"

$wsOverview.Range("A8").Value = "5030e9b0-8513-4b90-ab57-3936d9dca066.md"
$wsOverview.Range("B8").Value = "Ready for handoff"
$wsOverview.Range("C8").Value = "Ready for handoff"
$wsOverview.Range("D8").Value = "2016-28-12 12:28:57"

$wsOverview.Range("A9").Value = "b9bdd042-385e-47ff-88cc-5a87bc67efd0.md"
$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-27-12 12:27:09"

# zh-cn sheet (detail: File Name / File Extension / Status / Latest Handoff
# File / Latest Handoff Datetime / ... )
$wsZhCn = $excel.ActiveWorkbook.Worksheets.Item("zh-cn")

$wsZhCn.Range("A7").Value = "f843c243-6809-434c-80ed-ff15e90cb237.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "In Translation"
$wsZhCn.Range("D7").Value = "f843c243-6809-434c-80ed-ff15e90cb237.4c3edea4343f13c7713e0d09b8069b7f97e3d4d6.zh-cn.xlf"
$wsZhCn.Range("E7").Value = "2016-03-12 12:30:53"

$wsZhCn.Range("A8").Value = "5030e9b0-8513-4b90-ab57-3936d9dca066.md"
$wsZhCn.Range("B8").Value = ".md"
$wsZhCn.Range("C8").Value = "Ready for handoff"
$wsZhCn.Range("D8").Value = "5030e9b0-8513-4b90-ab57-3936d9dca066.013d84b898be606bdf9bc970d6e4e684e8c31c22.zh-cn.xlf"
$wsZhCn.Range("E8").Value = "2016-03-12 12:28:54"

$wsZhCn.Range("A9").Value = "b9bdd042-385e-47ff-88cc-5a87bc67efd0.md"
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "b9bdd042-385e-47ff-88cc-5a87bc67efd0.c5bf41561ef95c4176921bcb4581fd5bad49af23.zh-cn.xlf"
$wsZhCn.Range("E9").Value = "2016-03-12 12:27:05"

# de-de sheet (same layout as zh-cn)
$wsDeDe = $excel.ActiveWorkbook.Worksheets.Item("de-de")

$wsDeDe.Range("A7").Value = "f843c243-6809-434c-80ed-ff15e90cb237.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "In Translation"
$wsDeDe.Range("D7").Value = "f843c243-6809-434c-80ed-ff15e90cb237.4c3edea4343f13c7713e0d09b8069b7f97e3d4d6.de-de.xlf"
$wsDeDe.Range("E7").Value = "2016-03-12 12:30:56"

$wsDeDe.Range("A8").Value = "5030e9b0-8513-4b90-ab57-3936d9dca066.md"
$wsDeDe.Range("B8").Value = ".md"
$wsDeDe.Range("C8").Value = "Ready for handoff"
$wsDeDe.Range("D8").Value = "5030e9b0-8513-4b90-ab57-3936d9dca066.013d84b898be606bdf9bc970d6e4e684e8c31c22.de-de.xlf"
$wsDeDe.Range("E8").Value = "2016-03-12 12:28:57"

$wsDeDe.Range("A9").Value = "b9bdd042-385e-47ff-88cc-5a87bc67efd0.md"
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "b9bdd042-385e-47ff-88cc-5a87bc67efd0.c5bf41561ef95c4176921bcb4581fd5bad49af23.de-de.xlf"
$wsDeDe.Range("E9").Value = "2016-03-12 12:27:09"
